$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, shifting rows 37:45 down to 38:46
$ws.Rows.Item(37).Insert()

# Copy the static/repeated field values from row 38 (the row that got pushed down,
# which retains the same pattern) into the new row 37
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = "Vega Modelo de Temuco"
$ws.Range("C37").Value = "La Araucanía"
$ws.Range("D37").Value = 44505
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 7000
$ws.Range("M37").Value = 7000
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 280
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"

# Ensure date-formatted cell style carries over for the new row's date cell
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
